$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (number of passive trials per leg)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Row 2 (CON) updated measurements
$ws.Range("B2").Value = 13.093917949802517
$ws.Range("C2").Value = 18.945101293806943
$ws.Range("D2").Value = 26.223132047116678
$ws.Range("E2").Value = 27.891313829365004

# Row 3 (STR) updated measurements
$ws.Range("B3").Value = 11.493826343276282
$ws.Range("C3").Value = 20.104625531768875
$ws.Range("D3").Value = 17.00445050643782
$ws.Range("E3").Value = 30.781500941369586

# Selection now only covers the updated columns
$ws.Range("B1:E3").Select()
